# Refresh cryptos.xlsx price/volume snapshot (per-row Price [D] and
# Volume(1h) [E] updates), including the WrappedEther/WrappedBTC row swap
# (rows 17-18, columns B-E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.040.67"

$ws.Range("D3").Value = "3.132.58"

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'587.54"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "'146.09"
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "3.128.93"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("E9").Value = "  -1.45%  "

$ws.Range("E10").Value = "  -2.58%  "

$ws.Range("D11").Value = "'5.82"
$ws.Range("E11").Value = "  +2.29%  "

$ws.Range("E13").Value = "  -3.55%  "

$ws.Range("D14").Value = "'37.17"
$ws.Range("E14").Value = "  +3.20%  "

$ws.Range("D15").Value = "3.664.23"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "63.861.20"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.152.11"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").Value = "'7.07"
$ws.Range("E19").Value = "  -1.43%  "

$ws.Range("D20").Value = "'463.04"
$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("D21").Value = "'14.27"
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("E22").Value = "  -0.63%  "

$ws.Range("D23").Value = "'7.38"
$ws.Range("E23").Value = "  -2.23%  "

$ws.Range("E24").Value = "  -3.00%  "

$ws.Range("D25").Value = "'80.72"
$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  +5.68%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "'9.43"
$ws.Range("E28").Value = "  +8.97%  "

$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").Value = "'2.19"
$ws.Range("E31").Value = "  -1.07%  "

$ws.Range("D32").Value = "'7.11"
$ws.Range("E32").Value = "  +3.70%  "

$ws.Range("E33").Value = "  -0.51%  "

$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("D35").Value = "0.0₃0852"
$ws.Range("E35").Value = "  -2.22%  "

$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("E37").Value = "  -3.89%  "

$ws.Range("D38").Value = "'3.31"
$ws.Range("E38").Value = "  -2.77%  "

$ws.Range("D39").Value = "'6.01"
$ws.Range("E39").Value = "  -1.83%  "

$ws.Range("D40").Value = "'51.36"
$ws.Range("E40").Value = "  +1.03%  "

$ws.Range("D41").Value = "'438.18"
$ws.Range("E41").Value = "  -2.81%  "

$ws.Range("D42").Value = "'8.86"
$ws.Range("E42").Value = "  +1.72%  "

$ws.Range("D43").Value = "'0.285"
$ws.Range("E43").Value = "  +2.72%  "

$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").Value = "2.904.67"
$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").Value = "'39.64"
$ws.Range("E46").Value = "  +15.73%  "

$ws.Range("E47").Value = "  -3.64%  "

$ws.Range("D48").Value = "'126.52"
$ws.Range("E48").Value = "  -1.29%  "

$ws.Range("E50").Value = "  -0.92%  "

$ws.Range("E51").Value = "  +0.46%  "
